$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Uncertainty becomes a plain number (2) instead of text "2.0"
$ws.Range("D2").Value = 2

# Row 3: Value (B3) becomes text "0.6" instead of numeric 0.6
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.6"

# Row 3: Uncertainty (D3) becomes text "2.8" instead of numeric 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.8"
